# Append the 2026-02-15 bitcoin buy record as row 70.
# Column A stores the date as literal text (matching the existing rows,
# e.g. "02/08/2026" in A69) rather than a parsed date serial, so the cell
# is temporarily marked as Text before the value is typed in, then the
# style is reset back to Normal (matching the unstyled sibling cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCell = $ws.Cells.Item(70, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "02/15/2026"
$dateCell.Style = "Normal"

$ws.Cells.Item(70, 2).Value = 0.0007009300000000024
$ws.Cells.Item(70, 3).Value = 70620.46138701416
$ws.Cells.Item(70, 4).Value = 50
